$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new BOM line (10k-ohm resistor, designators R1/R2) was inserted as row 6
# of the BOM table, pushing the existing "MIC5504" regulator row down to row
# 7 and the trailing blank row down to row 8.
$ws.Rows(6).Insert()

# "No." column keeps its original text-typed look ("5", "6", ...) rather
# than becoming a literal number, so force Text formatting before typing the
# value, then drop back to the Normal style so no visible formatting change
# is left behind.
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "5"
$ws.Range("A6").Style = "Normal"

$ws.Range("B6").Value = 2
$ws.Range("C6").Value = "10kΩ"
$ws.Range("D6").Value = "R1,R2"
$ws.Range("E6").Value = "R0805"
$ws.Range("F6").Value = "10kΩ"
$ws.Range("G6").Value = "0805W8F1002T5E"
$ws.Range("H6").Value = "UNI-ROYAL(厚声)"
$ws.Range("I6").Value = "C17414"
$ws.Range("J6").Value = "LCSC"

# Renumber the regulator row that got pushed down to row 7.
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "6"
$ws.Range("A7").Style = "Normal"

# Rename the sheet/tab to reflect the new revision date (2024-12-1 -> 2024-12-2).
$ws.Name = "BOM_Board1_Schematic1_2024-12-2"
